# Daily auto push: insert a new timestamped reading for 2026/01/27 at row 735,
# shifting the existing rows 735:776 down to 736:777.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 735 (pushes old row 735..776 down to 736..777).
$ws.Rows.Item(735).Insert()

# Fill the new row with the new observation.
# Force text formatting on column A first so "2026/01/27" is stored as text
# (matching the rest of the column) rather than being auto-converted to a
# date serial number, then reset the style so no extra explicit cell style
# is left behind (matching the unstyled neighbouring date cells).
$ws.Range("A735").NumberFormat = "@"
$ws.Range("A735").Value = "2026/01/27"
$ws.Range("A735").Style = "Normal"

$ws.Range("B735").Value = "火"
$ws.Range("C735").Value = 19
$ws.Range("D735").Value = 201
